# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets to
# match the latest generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, identical on both sheets that
# carry the full data set (展览 = sheet 1, 全部类型 = sheet 4).
$updates = @{
    3  = 94
    4  = 277
    6  = 562
    7  = 58
    8  = 2028
    11 = 4380
    13 = 281
    14 = 101
    16 = 114
    19 = 71
    20 = 3184
    21 = 69
    22 = 469
    26 = 86
    29 = 55
    31 = 11
    32 = 560
    33 = 1804
    34 = 276
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
